$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Actor Input" text in the Flow of Events table (C8)
$ws.Range("C8").Value = "Indica que pretende ver as avaliações dos seus serviços"

# Widen column C (drop the bestFit auto-size, apply an explicit custom width)
$ws.Columns("C").ColumnWidth = 50.83

# Move the active cell selection to C9
$ws.Range("C9").Select()
